$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4, column A currently stores "79174445" as text (inline string).
# Convert it to a real number to match the updated data.
$ws.Cells.Item(4, 1).Value = 79174445

# Append new payment row 5 for phone 79174445 (Cash), paid 2025-08-23T09:41:10
# The phone number is kept as text, matching how it is stored elsewhere.
$ws.Cells.Item(5, 1).Value = "'79174445"
$ws.Cells.Item(5, 1).Style = "Normal"
$ws.Cells.Item(5, 2).Value = 3000
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 3000
$ws.Cells.Item(5, 7).Value = "Cash"
$ws.Cells.Item(5, 8).Value = "2025-08-23T09:41:10"
